$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("stats")

$ws.Range("D2").Value = 0.005102754570543766
$ws.Range("E2").Value = 0.08766797697171569
$ws.Range("G2").Value = 0.00647373590618372
$ws.Range("H2").Value = 0.05529298214241862
$ws.Range("I2").Value = 0.002181593794375658
$ws.Range("J2").Value = 0.01306318817660213
$ws.Range("K2").Value = 0.004212494008243084
$ws.Range("D3").Value = 0.002768766600638628
$ws.Range("E3").Value = 0.3978188671171665
$ws.Range("G3").Value = 0.0329343150369823
$ws.Range("H3").Value = 0.2477840362116694
$ws.Range("I3").Value = 0.02008050912991166
$ws.Range("J3").Value = 0.04348715161904693
$ws.Range("K3").Value = 0.02110619191080332
$ws.Range("D4").Value = 0.004544228781014681
$ws.Range("E4").Value = 0.08951445901766419
$ws.Range("G4").Value = 0.00708110723644495
$ws.Range("H4").Value = 0.05620518606156111
$ws.Range("I4").Value = 0.002070061396807432
$ws.Range("J4").Value = 0.01354377809911966
$ws.Range("K4").Value = 0.004198172129690647
$ws.Range("D5").Value = 0.00162486732006073
$ws.Range("E5").Value = 0.3988577299751341
$ws.Range("G5").Value = 0.03452201280742884
$ws.Range("H5").Value = 0.2477879533544183
$ws.Range("I5").Value = 0.02027707872912288
$ws.Range("J5").Value = 0.0441118972375989
$ws.Range("K5").Value = 0.02071312442421913
$ws.Range("C6").Value = 4805
$ws.Range("E6").Value = 1.298504695296288
$ws.Range("D7").Value = 0.002298991661518812
$ws.Range("E7").Value = 0.03687652293592691
$ws.Range("G7").Value = 0.002486546989530325
$ws.Range("H7").Value = 0.02200968703255057
$ws.Range("I7").Value = 0.001150201074779034
$ws.Range("J7").Value = 0.007398577872663736
$ws.Range("K7").Value = 0.00150013854727149
$ws.Range("D8").Value = 0.002187779173254967
$ws.Range("E8").Value = 0.3089819452725351
$ws.Range("G8").Value = 0.02513412712141871
$ws.Range("H8").Value = 0.1923001850955188
$ws.Range("I8").Value = 0.01554499007761478
$ws.Range("J8").Value = 0.03539645997807384
$ws.Range("K8").Value = 0.01592940557748079
$ws.Range("D9").Value = 0.002204176038503647
$ws.Range("E9").Value = 0.03687581699341536
$ws.Range("G9").Value = 0.00262652151286602
$ws.Range("H9").Value = 0.02256414387375116
$ws.Range("I9").Value = 0.001185137778520584
$ws.Range("J9").Value = 0.006533857434988022
$ws.Range("K9").Value = 0.00162496929988265
$ws.Range("D10").Value = 0.00132397934794426
$ws.Range("E10").Value = 0.3183064819313586
$ws.Range("G10").Value = 0.02731511415913701
$ws.Range("H10").Value = 0.1957804700359702
$ws.Range("I10").Value = 0.01657566009089351
$ws.Range("J10").Value = 0.03699039248749614
$ws.Range("K10").Value = 0.01647298689931631
$ws.Range("C11").Value = 4805
$ws.Range("E11").Value = 1.028471729718149
$ws.Range("D12").Value = 0.002532277256250381
$ws.Range("E12").Value = 0.04274594411253929
$ws.Range("G12").Value = 0.003110010176897049
$ws.Range("H12").Value = 0.02715241070836782
$ws.Range("I12").Value = 0.001401375513523817
$ws.Range("J12").Value = 0.006053993478417397
$ws.Range("K12").Value = 0.001905275508761406
$ws.Range("D13").Value = 0.001971160061657429
$ws.Range("E13").Value = 0.2753048120066524
$ws.Range("G13").Value = 0.02226255927234888
$ws.Range("H13").Value = 0.1715072728693485
$ws.Range("I13").Value = 0.01504562515765429
$ws.Range("J13").Value = 0.03031315561383963
$ws.Range("K13").Value = 0.01421913085505366
$ws.Range("D14").Value = 0.002352155745029449
$ws.Range("E14").Value = 0.04351937724277377
$ws.Range("G14").Value = 0.003388058859854937
$ws.Range("H14").Value = 0.02714584022760391
$ws.Range("I14").Value = 0.001327619422227144
$ws.Range("J14").Value = 0.006429675035178661
$ws.Range("K14").Value = 0.002142644021660089
$ws.Range("D15").Value = 0.001237661112099886
$ws.Range("E15").Value = 0.2934476262889802
$ws.Range("G15").Value = 0.02473073313012719
$ws.Range("H15").Value = 0.1807961696758866
$ws.Range("I15").Value = 0.01601315708830953
$ws.Range("J15").Value = 0.03334408765658736
$ws.Range("K15").Value = 0.01546579273417592
$ws.Range("C16").Value = 4805
$ws.Range("E16").Value = 1.0750294579193
$ws.Range("D17").Value = 0.003100324422121048
$ws.Range("E17").Value = 0.05447119008749723
$ws.Range("G17").Value = 0.003970513585954905
$ws.Range("H17").Value = 0.03474805178120732
$ws.Range("I17").Value = 0.001422002911567688
$ws.Range("J17").Value = 0.008151958230882883
$ws.Range("K17").Value = 0.002532115206122398
$ws.Range("D18").Value = 0.002204247750341892
$ws.Range("E18").Value = 0.316681609954685
$ws.Range("G18").Value = 0.02525182301178575
$ws.Range("H18").Value = 0.1976799173280597
$ws.Range("I18").Value = 0.01706172991544008
$ws.Range("J18").Value = 0.03536766255274415
$ws.Range("K18").Value = 0.01633666781708598
$ws.Range("D19").Value = 0.003161663189530373
$ws.Range("E19").Value = 0.05682329786941409
$ws.Range("G19").Value = 0.004329751711338758
$ws.Range("H19").Value = 0.0365743669681251
$ws.Range("I19").Value = 0.001428878866136074
$ws.Range("J19").Value = 0.008125264663249254
$ws.Range("K19").Value = 0.002553056925535202
$ws.Range("D20").Value = 0.001299806404858828
$ws.Range("E20").Value = 0.3178509171120822
$ws.Range("G20").Value = 0.02665189374238253
$ws.Range("H20").Value = 0.1965084415860474
$ws.Range("I20").Value = 0.01769693940877914
$ws.Range("J20").Value = 0.03630434768274426
$ws.Range("K20").Value = 0.01606337446719408
$ws.Range("C21").Value = 4805
$ws.Range("E21").Value = 1.042856383603066
$ws.Range("D22").Value = 0.002818004693835974
$ws.Range("E22").Value = 0.04500549519434571
$ws.Range("G22").Value = 0.00331448158249259
$ws.Range("H22").Value = 0.027995181735605
$ws.Range("I22").Value = 0.001283865422010422
$ws.Range("J22").Value = 0.007226055953651667
$ws.Range("K22").Value = 0.002052519004791975
$ws.Range("D23").Value = 0.002064861822873354
$ws.Range("E23").Value = 0.2926517301239073
$ws.Range("G23").Value = 0.02392783807590604
$ws.Range("H23").Value = 0.1824900843203068
$ws.Range("I23").Value = 0.01616172399371862
$ws.Range("J23").Value = 0.03163070045411587
$ws.Range("K23").Value = 0.01501797046512365
$ws.Range("D24").Value = 0.002560560591518879
$ws.Range("E24").Value = 0.04559829970821738
$ws.Range("G24").Value = 0.003526780288666487
$ws.Range("H24").Value = 0.02853199047967792
$ws.Range("I24").Value = 0.001291815191507339
$ws.Range("J24").Value = 0.007153824437409639
$ws.Range("K24").Value = 0.002053318545222282
$ws.Range("D25").Value = 0.001223662402480841
$ws.Range("E25").Value = 0.3002450447529554
$ws.Range("G25").Value = 0.02534763514995575
$ws.Range("H25").Value = 0.1827242709696293
$ws.Range("I25").Value = 0.01744035072624683
$ws.Range("J25").Value = 0.03507415438070893
$ws.Range("K25").Value = 0.01557184895500541
$ws.Range("C26").Value = 4805
$ws.Range("E26").Value = 1.037118395324796
$ws.Range("D27").Value = 0.003170470520853996
$ws.Range("E27").Value = 0.05335342586040497
$ws.Range("G27").Value = 0.003871057648211718
$ws.Range("H27").Value = 0.03343966268002987
$ws.Range("I27").Value = 0.001487807743251324
$ws.Range("J27").Value = 0.008378754742443562
$ws.Range("K27").Value = 0.002440508455038071
$ws.Range("D28").Value = 0.002239363081753254
$ws.Range("E28").Value = 0.3182877928949893
$ws.Range("G28").Value = 0.02590213250368834
$ws.Range("H28").Value = 0.1983522991649806
$ws.Range("I28").Value = 0.01677891565486789
$ws.Range("J28").Value = 0.03523902604356408
$ws.Range("K28").Value = 0.01652187332510948
$ws.Range("D29").Value = 0.002964556869119406
$ws.Range("E29").Value = 0.05446625016629696
$ws.Range("G29").Value = 0.004190443921834231
$ws.Range("H29").Value = 0.03420430552214384
$ws.Range("I29").Value = 0.001460702531039715
$ws.Range("J29").Value = 0.008357279933989048
$ws.Range("K29").Value = 0.002514432184398174
$ws.Range("D30").Value = 0.001341995317488909
$ws.Range("E30").Value = 0.3257415600121021
$ws.Range("G30").Value = 0.02771347779780626
$ws.Range("H30").Value = 0.2007194611243904
$ws.Range("I30").Value = 0.01760063720867038
$ws.Range("J30").Value = 0.03716497588902712
$ws.Range("K30").Value = 0.01685742549598217
$ws.Range("C31").Value = 4805
$ws.Range("E31").Value = 1.09639613237232
